$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add "Wins", "Losses", "Ties" in AD1:AF1 ---
# Copy the header style (bold font, thin border, centered/top-aligned)
# from an existing header cell so the new header cells match the rest
# of row 1 exactly (same style index, no new style entries created).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-39): season record for every player's team ---
# Wins = 93, Losses = 69, Ties = 0 for each row.
$lastRow = 39
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 69   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
